$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 117. This shifts the previous rows 117..194 down to 118..195,
# matching the dimension change from A1:T194 to A1:T195.
$ws.Rows.Item(117).Insert()

# The row that used to be 117 is now row 118. Copy its (unchanged) values into the
# new row 117, then overwrite the three fields (Fecha, Volumen, Origen) that differ
# for the newly inserted record.
for ($col = 1; $col -le 20; $col++) {
    $ws.Cells.Item(117, $col).Value = $ws.Cells.Item(118, $col).Value2
}

$ws.Cells.Item(117, 4).Value = 44606              # D117 Fecha
$ws.Cells.Item(117, 13).Value = 45                # M117 Volumen
$ws.Cells.Item(117, 18).Value = "Provincia de Linares"  # R117 Origen
